$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns at CE:CF, shifting the existing
# Outcome_Conservative_* / Outcome_Endoscopy_Surgery_* columns two places
# to the right (CE->CG, CF->CH, CG->CI, CH->CJ).
$ws.Range("CE1:CF1").EntireColumn.Insert()

# New header for the inserted "cases" column; the inserted rate column (CF1)
# was left carrying a copy of the (now shifted) "Outcome_Conservative_Rate"
# header text rather than being renamed yet.
$ws.Range("CE1").Value = "Outcome_Complication_Cases"
$ws.Range("CF1").Value = "Outcome_Conservative_Rate"

# Row 22 picked up explicit "UK" placeholder values in the two new cells.
$ws.Range("CE22").Value = "UK"
$ws.Range("CF22").Value = "UK"

# Row 26: first real data entered into the new columns.
$ws.Range("CE26").Value = 6
$ws.Range("CF26").Value = 0.158
$ws.Range("CF26").NumberFormat = "0.000"

# Reflect the author's navigation/view state at save time.
$ws.Range("CF26").Select()
$ws.Application.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("CF26").Select()
